# The "MASTER model" commit re-sorted the team-stats table; in this sheet
# that manifests as the Dallas Mavericks row (row 10) and the Milwaukee
# Bucks row (row 11) trading places. Column A just holds a 0-based row
# index (doesn't travel with the team), so we leave it alone. Columns B
# (LeagueID) and C (SeasonID) are identical for every team/row, so they're
# skipped too (swapping them is a no-op, and skipping avoids accidentally
# flipping their text storage to numeric). Every other populated column,
# D through CO, belongs to the team and needs to move with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rngRow10 = $ws.Range("D10:CO10")
$rngRow11 = $ws.Range("D11:CO11")

$row10Values = $rngRow10.Value2
$row11Values = $rngRow11.Value2

$rngRow10.Value2 = $row11Values
$rngRow11.Value2 = $row10Values
